$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Workbook window position (workbook.xml bookViews/workbookView@xWindow)
# ---------------------------------------------------------------------------
$excel.Left = 160

# ---------------------------------------------------------------------------
# 2) Re-style rows 85-88 so they match the "no-fill" look used elsewhere in
#    the table (row 15 / row 74 are good templates for this style family).
#    Values are left untouched here - only formats move.
# ---------------------------------------------------------------------------
$ws.Range("A15:E15").Copy() | Out-Null
$ws.Range("A85:E85").PasteSpecial(-4122) | Out-Null

$ws.Range("A74:E74").Copy() | Out-Null
$ws.Range("A86:E86").PasteSpecial(-4122) | Out-Null
$ws.Range("A87:E87").PasteSpecial(-4122) | Out-Null

$ws.Range("A75:E75").Copy() | Out-Null
$ws.Range("A88:E88").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Status-meeting updates to existing rows 84-88
# ---------------------------------------------------------------------------
$ws.Range("E84").Value = "In Progress"
$ws.Range("E85").Value = "Complete"
$ws.Range("E86").Value = "Complete - combined with 88"
$ws.Range("C88").Value = "Mervi Heiskanen, Shine Jacob"
$ws.Range("E88").Value = "Complete"

# ---------------------------------------------------------------------------
# 4) Fill in the three previously-blank rows (89-91) with new action items
# ---------------------------------------------------------------------------
$ws.Range("B89").Value = "Schedule meeting with Tabitha to discuss performance and usability."
$ws.Range("C89").Value = "Mike Hunter"
$ws.Range("D89").Value = 39960
$ws.Range("E89").Value = "Assigned"

$ws.Range("B90").Value = "Provide links to the two MAT KC videos related to data submission."
$ws.Range("C90").Value = "Jill Hadfield"
$ws.Range("D90").Value = 39960
$ws.Range("E90").Value = "Assigned"

$ws.Range("B91").Value = "Review the MAT KC videos related to data submission (links provided by Jill)."
$ws.Range("C91").Value = "Mike Hunter"
$ws.Range("D91").Value = 39960
$ws.Range("E91").Value = "Assigned"

# ---------------------------------------------------------------------------
# 5) Extend the blank-row template (row 92) down to row 102, matching the
#    existing formatting/numbering pattern.
# ---------------------------------------------------------------------------
$ws.Range("A92:E92").Copy() | Out-Null
$ws.Range("A93:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($r = 93; $r -le 102; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# ---------------------------------------------------------------------------
# 6) View state: scroll position + active selection
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A75"), $true)
$ws.Range("B100").Select()
